$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45938
$ws.Range("B2").Value = 4302.57804378519
$ws.Range("C2").Value = 5195.60278747806
$ws.Range("D2").Value = 7028
$ws.Range("E2").Value = 6278.998608
$ws.Range("F2").Value = 6.00097298720266

$ws.Range("A3").Value = 45939
$ws.Range("B3").Value = 4284.12570642337
$ws.Range("C3").Value = 4983.61608133537
$ws.Range("D3").Value = 3620
$ws.Range("E3").Value = 5922.178665
$ws.Range("F3").Value = 125.069543329667

$ws.Range("A4").Value = 45940
$ws.Range("B4").Value = 5266.41826360418
$ws.Range("C4").Value = 4589.48151536799
$ws.Range("D4").Value = 3620
$ws.Range("E4").Value = 7616.480829
$ws.Range("F4").Value = 138.314336698492

$ws.Range("A5").Value = 45941
$ws.Range("B5").Value = 1377.28021753348
$ws.Range("C5").Value = 3068.55822773336
$ws.Range("D5").Value = 3620
$ws.Range("E5").Value = 3281.365815
$ws.Range("F5").Value = 56.3601593833287

$ws.Range("A6").Value = 45942
$ws.Range("B6").Value = 1218.5054891698
$ws.Range("C6").Value = 3075.00596246482
$ws.Range("D6").Value = 3620
$ws.Range("E6").Value = 3041.057442
$ws.Range("F6").Value = 53.231579803959

$ws.Range("A7").Value = 45943
$ws.Range("B7").Value = 5841.28374289435
$ws.Range("C7").Value = 5551.83626593772
$ws.Range("D7").Value = 3620
$ws.Range("E7").Value = 8951.451256
$ws.Range("F7").Value = 210.083490793474

$ws.Range("A8").Value = 45944
$ws.Range("B8").Value = 5841.28374289435
$ws.Range("C8").Value = 6050.06053037344
$ws.Range("D8").Value = 3620
$ws.Range("E8").Value = 8951.451256
$ws.Range("F8").Value = 230.842835144962

$ws.Range("A9").Value = 45945
$ws.Range("B9").Value = 5841.28374289435
$ws.Range("C9").Value = 6519.0217646734
$ws.Range("D9").Value = 3620
$ws.Range("E9").Value = 8971.002178
$ws.Range("F9").Value = 251.197508324127

$ws.Range("A10").Value = 45946
$ws.Range("B10").Value = 5841.28374289435
$ws.Range("C10").Value = 6710.41756519948
$ws.Range("D10").Value = 3620
$ws.Range("E10").Value = 8971.002178
$ws.Range("F10").Value = 259.172333346047

$ws.Range("A11").Value = 45947
$ws.Range("B11").Value = 5841.28374289435
$ws.Range("C11").Value = 6143.88851450948
$ws.Range("D11").Value = 3620
$ws.Range("E11").Value = 8971.002178
$ws.Range("F11").Value = 235.566956233964

$ws.Range("A12").Value = 45948
$ws.Range("B12").Value = 1742.27770790123
$ws.Range("C12").Value = 4579.21969301854
$ws.Range("D12").Value = 3620
$ws.Range("E12").Value = 4473.847183
$ws.Range("F12").Value = 153.782882004888

$ws.Range("A13").Value = 45949
$ws.Range("B13").Value = 1636.94065696827
$ws.Range("C13").Value = 4453.21074769004
$ws.Range("D13").Value = 3620
$ws.Range("E13").Value = 4358.686082
$ws.Range("F13").Value = 148.123173863407

$ws.Range("A14").Value = 45950
$ws.Range("B14").Value = 6392.95297294923
$ws.Range("C14").Value = 6958.19855325215
$ws.Range("D14").Value = 3620
$ws.Range("E14").Value = 9928.516505
$ws.Range("F14").Value = 286.406753554288

$ws.Range("A15").Value = 45951
$ws.Range("B15").Value = 6392.95297294923
$ws.Range("C15").Value = 7130.84444460709
$ws.Range("D15").Value = 3620
$ws.Range("E15").Value = 9928.516505
$ws.Range("F15").Value = 293.600332360744
